$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.902067244052887
$ws.Range("B1").Value = 1.83492386341095
$ws.Range("C1").Value = 4.172513961791992
$ws.Range("D1").Value = 3.418251276016235
$ws.Range("E1").Value = 1.492251753807068
